$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1838.1538
$ws.Cells.Item(98, 9).Value = 1835.68
$ws.Cells.Item(98, 11).Value = 1835.68
$ws.Cells.Item(98, 13).Value = -337.6800000000001

$ws.Cells.Item(122, 8).Value = 1838.1538
$ws.Cells.Item(122, 9).Value = 1835.68
$ws.Cells.Item(122, 11).Value = 5507.04
$ws.Cells.Item(122, 13).Value = -3057.04

$ws.Cells.Item(132, 8).Value = 5119.222
$ws.Cells.Item(132, 9).Value = 5119.222
$ws.Cells.Item(132, 11).Value = 15357.666
$ws.Cells.Item(132, 13).Value = -12827.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2991.6323
$ws.Cells.Item(32, 9).Value = 2838.394
$ws.Cells.Item(32, 10).Value = 8048.5
$ws.Cells.Item(32, 11).Value = 2838.394
$ws.Cells.Item(32, 12).Value = 8048.5
$ws.Cells.Item(32, 13).Value = -2551.394
$ws.Cells.Item(32, 14).Value = -8622.5

$ws.Cells.Item(45, 8).Value = 16405.363
$ws.Cells.Item(45, 9).Value = 20042.04
$ws.Cells.Item(45, 11).Value = 20042.04
$ws.Cells.Item(45, 13).Value = -19665.04

$ws.Cells.Item(88, 8).Value = 6054.778
$ws.Cells.Item(88, 9).Value = 4597.6
$ws.Cells.Item(88, 10).Value = 7876.25
$ws.Cells.Item(88, 11).Value = 4597.6
$ws.Cells.Item(88, 12).Value = 7876.25
$ws.Cells.Item(88, 13).Value = -4191.6
$ws.Cells.Item(88, 14).Value = -8688.25

$ws.Cells.Item(91, 8).Value = 6054.778
$ws.Cells.Item(91, 9).Value = 4597.6
$ws.Cells.Item(91, 10).Value = 7876.25
$ws.Cells.Item(91, 11).Value = 4597.6
$ws.Cells.Item(91, 12).Value = 7876.25
$ws.Cells.Item(91, 13).Value = -3193.6
$ws.Cells.Item(91, 14).Value = -10684.25

$ws.Cells.Item(110, 8).Value = 3124.5757
$ws.Cells.Item(110, 9).Value = 2076.0715
$ws.Cells.Item(110, 11).Value = 2076.0715
$ws.Cells.Item(110, 13).Value = -31.07150000000001

$ws.Cells.Item(122, 8).Value = 4549.2
$ws.Cells.Item(122, 9).Value = 4290.3076
$ws.Cells.Item(122, 11).Value = 12870.9228
$ws.Cells.Item(122, 13).Value = -10420.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3944.6667
$ws.Cells.Item(134, 9).Value = 5667.3335
$ws.Cells.Item(134, 10).Value = 3083.3333
$ws.Cells.Item(134, 11).Value = 17002.0005
$ws.Cells.Item(134, 12).Value = 9249.999899999999
$ws.Cells.Item(134, 13).Value = -14467.0005
$ws.Cells.Item(134, 14).Value = -14319.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1539.725
$ws.Cells.Item(16, 9).Value = 1503.4642
$ws.Cells.Item(16, 10).Value = 1624.3334
$ws.Cells.Item(16, 11).Value = 1503.4642
$ws.Cells.Item(16, 12).Value = 1624.3334
$ws.Cells.Item(16, 13).Value = -1216.4642
$ws.Cells.Item(16, 14).Value = -2198.3334

$ws.Cells.Item(99, 8).Value = 7619.7
$ws.Cells.Item(99, 9).Value = 8593.200000000001
$ws.Cells.Item(99, 10).Value = 4699.2
$ws.Cells.Item(99, 11).Value = 8593.200000000001
$ws.Cells.Item(99, 12).Value = 4699.2
$ws.Cells.Item(99, 13).Value = -7095.200000000001
$ws.Cells.Item(99, 14).Value = -7695.2

$ws.Cells.Item(105, 8).Value = 2462.0625
$ws.Cells.Item(105, 9).Value = 2053.2222
$ws.Cells.Item(105, 10).Value = 2987.7144
$ws.Cells.Item(105, 11).Value = 2053.2222
$ws.Cells.Item(105, 12).Value = 2987.7144
$ws.Cells.Item(105, 13).Value = -306.2222000000002
$ws.Cells.Item(105, 14).Value = -6481.7144

$ws.Cells.Item(107, 8).Value = 4546132
$ws.Cells.Item(107, 9).Value = 6250745
$ws.Cells.Item(107, 11).Value = 6250745
$ws.Cells.Item(107, 13).Value = -6248825

$ws.Cells.Item(113, 8).Value = 1539.725
$ws.Cells.Item(113, 9).Value = 1503.4642
$ws.Cells.Item(113, 10).Value = 1624.3334
$ws.Cells.Item(113, 11).Value = 1503.4642
$ws.Cells.Item(113, 12).Value = 1624.3334
$ws.Cells.Item(113, 13).Value = 666.5358000000001
$ws.Cells.Item(113, 14).Value = -5964.3334

$ws.Cells.Item(122, 8).Value = 3083.5
$ws.Cells.Item(122, 9).Value = 2226.9167
$ws.Cells.Item(122, 10).Value = 4796.6665
$ws.Cells.Item(122, 11).Value = 6680.750100000001
$ws.Cells.Item(122, 12).Value = 14389.9995
$ws.Cells.Item(122, 13).Value = -4230.750100000001
$ws.Cells.Item(122, 14).Value = -19289.9995

$ws.Cells.Item(126, 8).Value = 7619.7
$ws.Cells.Item(126, 9).Value = 8593.200000000001
$ws.Cells.Item(126, 10).Value = 4699.2
$ws.Cells.Item(126, 11).Value = 25779.6
$ws.Cells.Item(126, 12).Value = 14097.6
$ws.Cells.Item(126, 13).Value = -23309.6
$ws.Cells.Item(126, 14).Value = -19037.6

$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 3136.074
$ws.Cells.Item(132, 9).Value = 2392.158
$ws.Cells.Item(132, 11).Value = 7176.474
$ws.Cells.Item(132, 13).Value = -4646.474

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 5113.5186
$ws.Cells.Item(113, 9).Value = 447.2857
$ws.Cells.Item(113, 11).Value = 1341.8571
$ws.Cells.Item(113, 13).Value = 828.1428999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 7177.971
$ws.Cells.Item(122, 9).Value = 6721.4443
$ws.Cells.Item(122, 10).Value = 8718.75
$ws.Cells.Item(122, 11).Value = 20164.3329
$ws.Cells.Item(122, 12).Value = 26156.25
$ws.Cells.Item(122, 13).Value = -17714.3329
$ws.Cells.Item(122, 14).Value = -31056.25

$ws.Cells.Item(132, 8).Value = 4943.5
$ws.Cells.Item(132, 9).Value = 1326.1666
$ws.Cells.Item(132, 11).Value = 3978.4998
$ws.Cells.Item(132, 13).Value = -1448.4998

$ws.Cells.Item(139, 8).Value = 99900.86
$ws.Cells.Item(139, 10).Value = 99900.86
$ws.Cells.Item(139, 12).Value = 99900.86
$ws.Cells.Item(139, 14).Value = -110180.86

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4337.6665
$ws.Cells.Item(7, 9).Value = 3536.25
$ws.Cells.Item(7, 10).Value = 8344.75
$ws.Cells.Item(7, 11).Value = 3536.25
$ws.Cells.Item(7, 12).Value = 8344.75
$ws.Cells.Item(7, 13).Value = -3424.25
$ws.Cells.Item(7, 14).Value = -8568.75

$ws.Cells.Item(122, 8).Value = 4081.4546
$ws.Cells.Item(122, 9).Value = 4189.7
$ws.Cells.Item(122, 10).Value = 2999
$ws.Cells.Item(122, 11).Value = 12569.1
$ws.Cells.Item(122, 12).Value = 8997
$ws.Cells.Item(122, 13).Value = -10119.1
$ws.Cells.Item(122, 14).Value = -13897

$ws.Cells.Item(126, 8).Value = 4337.6665
$ws.Cells.Item(126, 9).Value = 3536.25
$ws.Cells.Item(126, 10).Value = 8344.75
$ws.Cells.Item(126, 11).Value = 10608.75
$ws.Cells.Item(126, 12).Value = 25034.25
$ws.Cells.Item(126, 13).Value = -8138.75
$ws.Cells.Item(126, 14).Value = -29974.25

$ws.Cells.Item(131, 8).Value = 61428.57
$ws.Cells.Item(131, 9).Value = 60000
$ws.Cells.Item(131, 10).Value = 70000
$ws.Cells.Item(131, 11).Value = 60000
$ws.Cells.Item(131, 12).Value = 70000
$ws.Cells.Item(131, 13).Value = -54960
$ws.Cells.Item(131, 14).Value = -80080

$ws.Cells.Item(132, 8).Value = 5898.9585
$ws.Cells.Item(132, 9).Value = 4592
$ws.Cells.Item(132, 10).Value = 8512.875
$ws.Cells.Item(132, 11).Value = 13776
$ws.Cells.Item(132, 12).Value = 25538.625
$ws.Cells.Item(132, 13).Value = -11246
$ws.Cells.Item(132, 14).Value = -30598.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2202.4
$ws.Cells.Item(126, 9).Value = 2168.45
$ws.Cells.Item(126, 10).Value = 2338.2
$ws.Cells.Item(126, 11).Value = 6505.349999999999
$ws.Cells.Item(126, 12).Value = 7014.599999999999
$ws.Cells.Item(126, 13).Value = -4035.349999999999
$ws.Cells.Item(126, 14).Value = -11954.6
